$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (A1:D1) to snake_case field names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case the Spanish connector words (de/del/la/las/los/el/y) inside
#    state and municipality names, and normalize the grand-total label.
$ws.Range("B7").Value = 'Pabellón De Arteaga'
$ws.Range("B8").Value = 'Rincón De Romos'
$ws.Range("B9").Value = 'San Francisco De Los Romo'
$ws.Range("B28").Value = 'Amatenango De La Frontera'
$ws.Range("B32").Value = 'Bejucal De Ocampo'
$ws.Range("B34").Value = 'Benemérito De Las Américas'
$ws.Range("B41").Value = 'Chiapa De Corzo'
$ws.Range("B62").Value = 'Mazapa De Madero'
$ws.Range("B64").Value = 'Montecristo De Guerrero'
$ws.Range("B67").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B72").Value = 'Salto De Agua'
$ws.Range("B74").Value = 'San Cristóbal De Las Casas'
$ws.Range("B100").Value = 'Coyame Del Sotol'
$ws.Range("B107").Value = 'Hidalgo Del Parral'
$ws.Range("A120").Value = 'Ciudad De México'
$ws.Range("B123").Value = 'Cuajimalpa De Morelos'
$ws.Range("A138").Value = 'Coahuila De Zaragoza'
$ws.Range("B148").Value = 'San Juan De Sabinas'
$ws.Range("B170").Value = 'Nombre De Dios'
$ws.Range("B178").Value = 'San Juan De Guadalupe'
$ws.Range("B179").Value = 'San Juan Del Río'
$ws.Range("B180").Value = 'San Luis Del Cordero'
$ws.Range("A189").Value = 'Estado De México'
$ws.Range("B189").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B191").Value = 'Almoloya De Juárez'
$ws.Range("B195").Value = 'Atizapán De Zaragoza'
$ws.Range("B205").Value = 'Ecatepec De Morelos'
$ws.Range("B208").Value = 'Ixtapan De La Sal'
$ws.Range("B209").Value = 'Ixtapan Del Oro'
$ws.Range("B217").Value = 'Naucalpan De Juárez'
$ws.Range("B224").Value = 'San Felipe Del Progreso'
$ws.Range("B225").Value = 'San Martín De Las Pirámides'
$ws.Range("B231").Value = 'Tenango Del Valle'
$ws.Range("B237").Value = 'Tlalnepantla De Baz'
$ws.Range("B241").Value = 'Valle De Bravo'
$ws.Range("B242").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B245").Value = 'Villa De Allende'
$ws.Range("B246").Value = 'Villa Del Carbón'
$ws.Range("B254").Value = 'Apaseo El Alto'
$ws.Range("B255").Value = 'Apaseo El Grande'
$ws.Range("B263").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B272").Value = 'Purísima Del Rincón'
$ws.Range("B277").Value = 'San Diego De La Unión'
$ws.Range("B279").Value = 'San Francisco Del Rincón'
$ws.Range("B281").Value = 'San Luis De La Paz'
$ws.Range("B282").Value = 'San Miguel De Allende'
$ws.Range("B283").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B285").Value = 'Silao De La Victoria'
$ws.Range("B290").Value = 'Valle De Santiago'
$ws.Range("B296").Value = 'Acapulco De Juárez'
$ws.Range("B298").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B299").Value = 'Alcozauca De Guerrero'
$ws.Range("B302").Value = 'Atenango Del Río'
$ws.Range("B304").Value = 'Atoyac De Álvarez'
$ws.Range("B305").Value = 'Ayutla De Los Libres'
$ws.Range("B307").Value = 'Chilapa De Álvarez'
$ws.Range("B308").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B313").Value = 'Coyuca De Benítez'
$ws.Range("B314").Value = 'Coyuca De Catalán'
$ws.Range("B316").Value = 'Cuetzala Del Progreso'
$ws.Range("B317").Value = 'Cutzamala De Pinzón'
$ws.Range("B324").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B325").Value = 'Iguala De La Independencia'
$ws.Range("B327").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B328").Value = 'José Joaquín De Herrera'
$ws.Range("B342").Value = 'Taxco De Alarcón'
$ws.Range("B345").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B347").Value = 'Tixtla De Guerrero'
$ws.Range("B350").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B351").Value = 'Tlapa De Comonfort'
$ws.Range("B353").Value = 'Técpan De Galeana'
$ws.Range("B355").Value = 'Zihuatanejo De Azueta'
$ws.Range("B361").Value = 'Agua Blanca De Iturbide'
$ws.Range("B367").Value = 'Atotonilco El Grande'
$ws.Range("B370").Value = 'Cuautepec De Hinojosa'
$ws.Range("B376").Value = 'Huejutla De Reyes'
$ws.Range("B379").Value = 'Jacala De Ledezma'
$ws.Range("B384").Value = 'Mineral De La Reforma'
$ws.Range("B385").Value = 'Mineral Del Monte'
$ws.Range("B387").Value = 'Pachuca De Soto'
$ws.Range("B389").Value = 'Progreso De Obregón'
$ws.Range("B392").Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range("B397").Value = 'Tepehuacán De Guerrero'
$ws.Range("B398").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B400").Value = 'Tezontepec De Aldama'
$ws.Range("B404").Value = 'Tula De Allende'
$ws.Range("B405").Value = 'Tulancingo De Bravo'
$ws.Range("B406").Value = 'Villa De Tezontepec'
$ws.Range("B407").Value = 'Zacualtipán De Ángeles'
$ws.Range("B414").Value = 'Atotonilco El Alto'
$ws.Range("B415").Value = 'Autlán De Navarro'
$ws.Range("B424").Value = 'Encarnación De Díaz'
$ws.Range("B426").Value = 'Huejuquilla El Alto'
$ws.Range("B432").Value = 'Lagos De Moreno'
$ws.Range("B436").Value = 'Ojuelos De Jalisco'
$ws.Range("B440").Value = 'San Diego De Alejandría'
$ws.Range("B443").Value = 'San Martín De Bolaños'
$ws.Range("B444").Value = 'San Miguel El Alto'
$ws.Range("B446").Value = 'Santa María De Los Ángeles'
$ws.Range("B448").Value = 'Talpa De Allende'
$ws.Range("B449").Value = 'Tamazula De Gordiano'
$ws.Range("B452").Value = 'Tepatitlán De Morelos'
$ws.Range("B454").Value = 'Tizapán El Alto'
$ws.Range("B459").Value = 'Unión De San Antonio'
$ws.Range("B460").Value = 'Unión De Tula'
$ws.Range("B461").Value = 'Valle De Juárez'
$ws.Range("B466").Value = 'Zapotlán El Grande'
$ws.Range("A468").Value = 'Michoacán De Ocampo'
$ws.Range("B524").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B556").Value = 'Puente De Ixtla'
$ws.Range("B559").Value = 'Tetela Del Volcán'
$ws.Range("B560").Value = 'Tlaltizapán De Zapata'
$ws.Range("B567").Value = 'Bahía De Banderas'
$ws.Range("B570").Value = 'Ixtlán Del Río'
$ws.Range("B574").Value = 'Santa María Del Oro'
$ws.Range("B593").Value = 'Mier Y Noriega'
$ws.Range("B596").Value = 'San Nicolás De Los Garza'
$ws.Range("B598").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B603").Value = 'Cuilápam De Guerrero'
$ws.Range("B604").Value = 'Fresnillo De Trujano'
$ws.Range("B605").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B606").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B607").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B608").Value = 'Huajuapan De León'
$ws.Range("B610").Value = 'Ixtlán De Juárez'
$ws.Range("B614").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B615").Value = 'Oaxaca De Juárez'
$ws.Range("B616").Value = 'Ocotlán De Morelos'
$ws.Range("B618").Value = 'Putla Villa De Guerrero'
$ws.Range("B619").Value = 'Reforma De Pineda'
$ws.Range("B624").Value = 'San Antonino El Alto'
$ws.Range("B625").Value = 'San Antonio De La Cal'
$ws.Range("B693").Value = 'Santo Domingo De Morelos'
$ws.Range("B698").Value = 'Tataltepec De Valdés'
$ws.Range("B699").Value = 'Teotitlán De Flores Magón'
$ws.Range("B700").Value = 'Tlacolula De Matamoros'
$ws.Range("B701").Value = 'Totontepec Villa De Morelos'
$ws.Range("B702").Value = 'Villa Sola De Vega'
$ws.Range("B703").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B704").Value = 'Villa De Zaachila'
$ws.Range("B705").Value = 'Zimatlán De Álvarez'
$ws.Range("B728").Value = 'Cuetzalan Del Progreso'
$ws.Range("B733").Value = 'Huehuetlán El Chico'
$ws.Range("B734").Value = 'Huehuetlán El Grande'
$ws.Range("B739").Value = 'Ixcamilpa De Guerrero'
$ws.Range("B742").Value = 'Izúcar De Matamoros'
$ws.Range("B762").Value = 'Tepexi De Rodríguez'
$ws.Range("B763").Value = 'Teteles De Avila Castillo'
$ws.Range("B766").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B783").Value = 'Amealco De Bonfil'
$ws.Range("B785").Value = 'Cadereyta De Montes'
$ws.Range("B790").Value = 'Jalpan De Serra'
$ws.Range("B791").Value = 'Landa De Matamoros'
$ws.Range("B794").Value = 'Pinal De Amoles'
$ws.Range("B797").Value = 'San Juan Del Río'
$ws.Range("B806").Value = 'Axtla De Terrazas'
$ws.Range("B812").Value = 'Ciudad Del Maíz'
$ws.Range("B818").Value = 'Mexquitic De Carmona'
$ws.Range("B823").Value = 'San Ciro De Acosta'
$ws.Range("B829").Value = 'Santa María Del Río'
$ws.Range("B831").Value = 'Soledad De Graciano Sánchez'
$ws.Range("B836").Value = 'Tanquián De Escobedo'
$ws.Range("B841").Value = 'Villa De Arista'
$ws.Range("B842").Value = 'Villa De Arriaga'
$ws.Range("B843").Value = 'Villa De Guadalupe'
$ws.Range("B844").Value = 'Villa De Ramos'
$ws.Range("B845").Value = 'Villa De Reyes'
$ws.Range("B846").Value = 'Villa De La Paz'
$ws.Range("B907").Value = 'Soto La Marina'
$ws.Range("B918").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("A927").Value = 'Veracruz De Ignacio De La Llave'
$ws.Range("B933").Value = 'Amatlán De Los Reyes'
$ws.Range("B938").Value = 'Boca Del Río'
$ws.Range("B940").Value = 'Castillo De Teayo'
$ws.Range("B942").Value = 'Cazones De Herrera'
$ws.Range("B952").Value = 'Cosamaloapan De Carpio'
$ws.Range("B966").Value = 'Hueyapan De Ocampo'
$ws.Range("B967").Value = 'Ignacio De La Llave'
$ws.Range("B970").Value = 'Ixhuacán De Los Reyes'
$ws.Range("B971").Value = 'Ixhuatlán De Madero'
$ws.Range("B972").Value = 'Ixhuatlán Del Sureste'
$ws.Range("B979").Value = 'Juchique De Ferrer'
$ws.Range("B983").Value = 'Lerdo De Tejada'
$ws.Range("B985").Value = 'Martínez De La Torre'
$ws.Range("B995").Value = 'Paso De Ovejas'
$ws.Range("B998").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1004").Value = 'Sayula De Alemán'
$ws.Range("B1005").Value = 'Soledad De Doblado'
$ws.Range("B1040").Value = 'Concepción Del Oro'
$ws.Range("B1047").Value = 'Jiménez Del Teul'
$ws.Range("B1055").Value = 'Moyahua De Estrada'
$ws.Range("B1056").Value = 'Nochistlán De Mejía'
$ws.Range("B1067").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1072").Value = 'Villa De Cos'
$ws.Range("A1076").Value = 'Total'

# 3. Drop the footnote rows (sample size / source / author / date) that
#    trailed the data table, and shrink the used range back to the table.
$ws.Range("A1078:A1082").EntireRow.Delete() | Out-Null

